# Verify_24V_Load_On_Changing_CPU_Of_Panel.xlsx
# "Updated 24V test cases and test data with new loading details method"
#
# All the data edits land on the workbook's first/active sheet ("Add Panels").
# We add a new "Loading Details Name" column (K) with a header and two
# "Main Processor 24V (A)" values for the existing data rows, matching the
# formatting of neighbouring cells, then move the active selection to G7.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- K7: header cell, same style as the other header cells in row 7 (e.g. A7) ---
$ws.Range("A7").Copy()
$ws.Range("K7").PasteSpecial(-4122)   # xlPasteFormats
$ws.Range("K7").Value = "Loading Details Name"

# --- K8: data cell. Base it on C4's fill/border, then add left align + wrap ---
$ws.Range("C4").Copy()
$ws.Range("K8").PasteSpecial(-4122)   # xlPasteFormats
$ws.Range("K8").HorizontalAlignment = -4131   # xlLeft
$ws.Range("K8").WrapText = $true
$ws.Range("K8").Value = "Main Processor 24V (A)"

# --- K9: same style as K8 ---
$ws.Range("K8").Copy()
$ws.Range("K9").PasteSpecial(-4122)   # xlPasteFormats
$ws.Range("K9").Value = "Main Processor 24V (A)"

$excel.CutCopyMode = $false

# Size the new column to fit its contents (like the other bestFit columns)
$ws.Columns.Item(11).EntireColumn.AutoFit()

# Move the active selection (was C13) to G7
$ws.Range("G7").Select()
